# The document contains the field text "{m:'A sample table'.bothMergeAll()}"
# spread across several runs. Two of those runs need to be split in two
# (without altering their run formatting):
#   "{m"  -> "{" + "m"
#   "()}" -> "()" + "}"
#
# Word COM has no direct "split run" API, but inserting (and immediately
# removing) a temporary bookmark at the split point forces Word to break
# the run there without touching the run's formatting (no stray <w:rPr/>
# and no leftover bookmark markup once it is deleted).

$d = $word.ActiveDocument

function Split-RunAfter($searchText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find '$searchText' to split after"
    }
    $d.Bookmarks.Add("tmpSplitMark", $rng)
    $d.Bookmarks("tmpSplitMark").Delete()
}

# "{m" -> "{" | "m"
Split-RunAfter "{"

# "()}" -> "()" | "}"
Split-RunAfter "()"
